$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: I0 (I) and IF (J), matching the styling of the existing
# header row (bold, bordered, centered) by copying format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row, I-value, J-value for rows 2..68
$rowsData = @(
    @(2, 4, 5),
    @(3, 3, 4),
    @(4, 9, 9),
    @(5, 9, 9),
    @(6, 6, 6),
    @(7, 10, 10),
    @(8, 9, 9),
    @(9, 9, 10),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 8, 9),
    @(16, 8, 9),
    @(17, 10, 11),
    @(18, 7, 7),
    @(19, 4, 4),
    @(20, 7, 8),
    @(21, 8, 8),
    @(22, 6, 6),
    @(23, 10, 10),
    @(24, 7, 7),
    @(25, 6, 6),
    @(26, 5, 6),
    @(27, 6, 6),
    @(28, 5, 6),
    @(29, 5, 6),
    @(30, 7, 8),
    @(31, 3, 3),
    @(32, 7, 7),
    @(33, 1, 2),
    @(34, 6, 6),
    @(35, 1, 1),
    @(36, 4, 5),
    @(37, 7, 7),
    @(38, 8, 8),
    @(39, 7, 7),
    @(40, 9, 10),
    @(41, 5, 6),
    @(42, 10, 11),
    @(43, 6, 7),
    @(44, 11, 11),
    @(45, 7, 7),
    @(46, 6, 7),
    @(47, 4, 5),
    @(48, 6, 7),
    @(49, 8, 9),
    @(50, 7, 7),
    @(51, 1, 2),
    @(52, 8, 8),
    @(53, 7, 8),
    @(54, 9, 9),
    @(55, 3, 3),
    @(56, 1, 2),
    @(57, 5, 6),
    @(58, 3, 4),
    @(59, 10, 10),
    @(60, 8, 8),
    @(61, 7, 7),
    @(62, 7, 7),
    @(63, 7, 7),
    @(64, 4, 5),
    @(65, 6, 6),
    @(66, 7, 7),
    @(67, 4, 4),
    @(68, 3, 3)
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "Added I0/IF columns (I1:J68)"
